$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Create the "NormalWeb" ("Normal (Web)") paragraph style.
#    Reference it from a paragraph first so the engine mints it
#    without the w:customStyle="1" flag (matches a built-in style
#    id), then fetch the minted style object and fill in the rest
#    of its definition explicitly.
# ---------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Style = "NormalWeb"

$nw = $d.Styles("NormalWeb")
$nw.NameLocal = "Normal (Web)"
$nw.BaseStyle = "Normal"
$nw.Priority = 99
$nw.UnhideWhenUsed = $true
$nw.QuickStyle = $false
$nw.Font.Name = "Times New Roman"
$nw.Font.NameFarEast = "Times New Roman"
$nw.Font.NameBi = "Times New Roman"
$nw.Font.Size = 12
$nw.Font.SizeBi = 12
$nw.ParagraphFormat.SpaceBefore = 5
$nw.ParagraphFormat.SpaceBeforeAuto = $true
$nw.ParagraphFormat.SpaceAfter = 5
$nw.ParagraphFormat.SpaceAfterAuto = $true
$nw.ParagraphFormat.LineSpacingRule = 0

# ---------------------------------------------------------------
# 2. Rewrite the first body paragraph.
# ---------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "During this project it has given me the opportunity to do some research on what some of the risks could be in nanotechnology in medicine. I've never heard of Nanotechnology and was a bit intimidated by the name once we chose this as our group topic, however, I was determined to explore what the world of nanotechnology had to teach me."
$p1.Style = "NormalWeb"
$p1.Range.Font.Color = 0
$p1.Range.Font.Size = 13.5
$p1.Range.Font.SizeBi = 13.5

# ---------------------------------------------------------------
# 3. Rewrite the second body paragraph.
# ---------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "I found that an ethical issue in nanomedicine is protecting patients private information where medical records of body organs are stored electronically of the patients results. In addition, it has been revealed that when clinicians have a consultation remotely with their patients to discuss test results or diagnosis that has recently been discovered, this data is already stored in the system which is allocated against the patients file. The quantity of the storage is quite large, therefore it is recommended to have a highly sourced secure system to protect patient information from a cyber-attack."
$p2.Style = "NormalWeb"
$p2.Range.Font.Color = 0
$p2.Range.Font.Size = 13.5
$p2.Range.Font.SizeBi = 13.5

# ---------------------------------------------------------------
# 4. Remove the third body paragraph entirely (the old
#    "One recommendation ..." paragraph). The trailing empty
#    paragraph stays untouched.
# ---------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.Delete()

# ---------------------------------------------------------------
# 5. Update the header: name + student id, tagging the edited
#    runs with English (US) language.
# ---------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$hdrRange = $hdr.Range
$hdrRange.Find.Execute("Von Macatuhay", $false, $false, $false, $false, $false, $true, 1, $false, "Kika Kalolo", 2)

$hdrRange2 = $hdr.Range
$hdrRange2.Find.Execute("Kika Kalolo", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($hdrRange2.Find.Found) {
    $hdrRange2.LanguageID = "en-US"
}

$hdrRange3 = $hdr.Range
$hdrRange3.Find.Execute("19078493", $false, $false, $false, $false, $false, $true, 1, $false, "19088935", 2)

$hdrRange4 = $hdr.Range
$hdrRange4.Find.Execute("19088935", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($hdrRange4.Find.Found) {
    $hdrRange4.LanguageID = "en-US"
}

Write-Host "done"
